$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they stay text like the originals.
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D16", "D20", "D22", "D23", "D24", "D25", "D26", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D47", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.408.42'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '2.220.88'
$ws.Range("E3").Value = '  -2.18%  '
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '109.31'
$ws.Range("E5").Value = '  -7.84%  '
$ws.Range("D6").Value = '287.94'
$ws.Range("E6").Value = '  +7.30%  '
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  -3.39%  '
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("D9").Value = '0.592'
$ws.Range("E9").Value = '  -4.68%  '
$ws.Range("D10").Value = '43.08'
$ws.Range("E10").Value = '  -9.00%  '
$ws.Range("D11").Value = '0.0908'
$ws.Range("E11").Value = '  -4.08%  '
$ws.Range("D12").Value = '54.16'
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = '8.55'
$ws.Range("E13").Value = '  -9.25%  '
$ws.Range("E14").Value = '  +9.67%  '
$ws.Range("E15").Value = '  -2.91%  '
$ws.Range("D16").Value = '14.83'
$ws.Range("E16").Value = '  -6.84%  '
$ws.Range("D17").Value = '2.555.66'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").Value = '2.231.70'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").Value = '42.377.81'
$ws.Range("E19").Value = '  -2.85%  '
$ws.Range("D20").Value = '7.10'
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("E21").Value = '  -5.15%  '
$ws.Range("D22").Value = '72.82'
$ws.Range("E22").Value = '  +0.60%  '
$ws.Range("D23").Value = '3.36'
$ws.Range("E23").Value = '  +12.81%  '
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("D25").Value = '229.30'
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("D26").Value = '8.87'
$ws.Range("E26").Value = '  -8.32%  '
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("D28").Value = '11.35'
$ws.Range("E28").Value = '  -8.32%  '
$ws.Range("E29").Value = '  -2.67%  '
$ws.Range("D30").Value = '172.76'
$ws.Range("E30").Value = '  -1.08%  '
$ws.Range("B31").Value = 'WEMIXToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").Value = '3.15'
$ws.Range("E31").Value = '  -6.07%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '36.87'
$ws.Range("E32").Value = '  -11.86%  '
$ws.Range("D33").Value = '20.71'
$ws.Range("E33").Value = '  -3.87%  '
$ws.Range("D34").Value = '0.0870'
$ws.Range("E34").Value = '  -5.68%  '
$ws.Range("D35").Value = '5.56'
$ws.Range("E35").Value = '  -2.98%  '
$ws.Range("D36").Value = '4.91'
$ws.Range("E36").Value = '  +5.41%  '
$ws.Range("E37").Value = '  -3.94%  '
$ws.Range("D38").Value = '4.11'
$ws.Range("E38").Value = '  -4.17%  '
$ws.Range("D39").Value = '0.0365'
$ws.Range("E39").Value = '  -3.99%  '
$ws.Range("D40").Value = '0.103'
$ws.Range("E40").Value = '  -5.07%  '
$ws.Range("D41").Value = '73.27'
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("D42").Value = '2.36'
$ws.Range("E42").Value = '  -7.93%  '
$ws.Range("D43").Value = '0.228'
$ws.Range("E43").Value = '  -5.15%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '12.21'
$ws.Range("E45").Value = '  -11.66%  '
$ws.Range("E46").Value = '  -6.79%  '
$ws.Range("D47").Value = '5.29'
$ws.Range("E47").Value = '  -7.75%  '
$ws.Range("E48").Value = '  +8.39%  '
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '8.38'
$ws.Range("E50").Value = '  -2.41%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '100.68'
$ws.Range("E51").Value = '  -2.62%  '
